$d = $word.ActiveDocument

# Helper: replace the text-run contents of a paragraph (leaving its <w:pPr>,
# paraId, etc. untouched) with the raw run-level OOXML supplied in $xmlInner.
function Set-ParagraphRuns($paragraph, $xmlInner) {
    $full = $paragraph.Range
    $full.MoveEnd(1, -1) | Out-Null   # exclude the paragraph mark
    $startPos = $full.Start
    $full.Delete()

    $ins = $d.Range($startPos, $startPos)
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
        $xmlInner +
        '</w:p></w:body></w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
    $ins.InsertXML($xml)
}

$apos = [char]0x2019
$styleRunOpen = '<w:r><w:rPr><w:rStyle w:val="SubtleReference"/></w:rPr>'

# 1) "Bag:" -> "needs opened" becomes "Doesn't" + " need" + " opened" (3 runs)
$p = $d.Paragraphs.Item(32)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "needs opened") {
    $inner = $styleRunOpen + '<w:t>Doesn' + $apos + 't</w:t></w:r>' +
             $styleRunOpen + '<w:t xml:space="preserve"> need</w:t></w:r>' +
             $styleRunOpen + '<w:t xml:space="preserve"> opened</w:t></w:r>'
    Set-ParagraphRuns $p $inner
}

# 2) "Bag:" -> "Cant be locked" becomes proofErr-wrapped "Cant" + " be locked"
$p = $d.Paragraphs.Item(33)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Cant be locked") {
    $inner = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
             $styleRunOpen + '<w:t>Cant</w:t></w:r>' +
             '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
             $styleRunOpen + '<w:t xml:space="preserve"> be locked</w:t></w:r>'
    Set-ParagraphRuns $p $inner
}

# 3) "Barrel:" -> "Cant be locked" becomes proofErr-wrapped "Cant" + " be locked"
$p = $d.Paragraphs.Item(44)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Cant be locked") {
    $inner = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
             $styleRunOpen + '<w:t>Cant</w:t></w:r>' +
             '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
             $styleRunOpen + '<w:t xml:space="preserve"> be locked</w:t></w:r>'
    Set-ParagraphRuns $p $inner
}

# 4) "Barrel:" -> "Cant be picked up" becomes proofErr-wrapped "Cant" + " be picked up"
$p = $d.Paragraphs.Item(45)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Cant be picked up") {
    $inner = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
             $styleRunOpen + '<w:t>Cant</w:t></w:r>' +
             '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
             $styleRunOpen + '<w:t xml:space="preserve"> be picked up</w:t></w:r>'
    Set-ParagraphRuns $p $inner
}

# 5) "Chest:" -> "Cant be picked up" becomes proofErr-wrapped "Cant" + " be picked up"
$p = $d.Paragraphs.Item(50)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Cant be picked up") {
    $inner = '<w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/>' +
             $styleRunOpen + '<w:t>Cant</w:t></w:r>' +
             '<w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/>' +
             $styleRunOpen + '<w:t xml:space="preserve"> be picked up</w:t></w:r>'
    Set-ParagraphRuns $p $inner
}

# 6) Remove the paragraph "Will also need to update entity class to have Type
#    variable in order to differentiate these." entirely (text + paragraph mark).
$p = $d.Paragraphs.Item(54)
if ($p.Range.Text.TrimEnd([char]13,[char]7) -eq "Will also need to update entity class to have Type variable in order to differentiate these.") {
    $p.Range.Delete()
}
